# Auto-generated edit script: updates currentAveragePrice / Leve price / profit
# columns (H-N) across 8 item sheets to match refreshed market data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value2 = 233.33333
$ws.Range("I9").Value2 = 0
$ws.Range("J9").Value2 = 233.33333
$ws.Range("K9").Value2 = 0
$ws.Range("L9").Value2 = 233.33333
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value2 = -571.3333299999999
$ws.Range("H64").Value2 = 4631.3335
$ws.Range("I64").Value2 = 4324.5
$ws.Range("K64").Value2 = 4324.5
$ws.Range("M64").Value2 = -4076.5
$ws.Range("H67").Value2 = 4631.3335
$ws.Range("I67").Value2 = 4324.5
$ws.Range("K67").Value2 = 4324.5
$ws.Range("M67").Value2 = -3466.5
$ws.Range("H70").Value2 = 26517308
$ws.Range("J70").Value2 = 19609934
$ws.Range("L70").Value2 = 58829802
$ws.Range("N70").Value2 = -58830342
$ws.Range("H73").Value2 = 26517308
$ws.Range("J73").Value2 = 19609934
$ws.Range("L73").Value2 = 58829802
$ws.Range("N73").Value2 = -58831674
$ws.Range("H87").Value2 = 49999
$ws.Range("J87").Value2 = 49999
$ws.Range("L87").Value2 = 49999
$ws.Range("N87").Value2 = -52495
$ws.Range("H90").Value2 = 49999
$ws.Range("J90").Value2 = 49999
$ws.Range("L90").Value2 = 149997
$ws.Range("N90").Value2 = -162477
$ws.Range("H107").Value2 = 18549314
$ws.Range("I107").Value2 = 9000976
$ws.Range("K107").Value2 = 9000976
$ws.Range("M107").Value2 = -8999056

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 2505349.2
$ws.Range("I32").Value2 = 2556465.5
$ws.Range("K32").Value2 = 2556465.5
$ws.Range("M32").Value2 = -2556178.5
$ws.Range("H61").Value2 = 4674.0586
$ws.Range("I61").Value2 = 1919.1842
$ws.Range("K61").Value2 = 1919.1842
$ws.Range("M61").Value2 = -1707.1842
$ws.Range("H63").Value2 = 2895.6
$ws.Range("J63").Value2 = 3500.5
$ws.Range("L63").Value2 = 3500.5
$ws.Range("N63").Value2 = -4872.5
$ws.Range("H66").Value2 = 2895.6
$ws.Range("J66").Value2 = 3500.5
$ws.Range("L66").Value2 = 17502.5
$ws.Range("N66").Value2 = -24366.5
$ws.Range("H132").Value2 = 7125.514
$ws.Range("I132").Value2 = 6395.9414
$ws.Range("J132").Value2 = 7814.5557
$ws.Range("K132").Value2 = 19187.8242
$ws.Range("L132").Value2 = 23443.6671
$ws.Range("M132").Value2 = -16657.8242
$ws.Range("N132").Value2 = -28503.6671
$ws.Range("H136").Value2 = 4674.0586
$ws.Range("I136").Value2 = 1919.1842
$ws.Range("K136").Value2 = 5757.5526
$ws.Range("M136").Value2 = -3207.5526

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value2 = 7248107.5
$ws.Range("I20").Value2 = 18520548
$ws.Range("J20").Value2 = 1539.0714
$ws.Range("K20").Value2 = 18520548
$ws.Range("L20").Value2 = 1539.0714
$ws.Range("M20").Value2 = -18520301
$ws.Range("N20").Value2 = -2033.0714
$ws.Range("H99").Value2 = 4136487.2
$ws.Range("I99").Value2 = 4402.1055
$ws.Range("J99").Value2 = 30306360
$ws.Range("K99").Value2 = 4402.1055
$ws.Range("L99").Value2 = 30306360
$ws.Range("M99").Value2 = -2904.1055
$ws.Range("N99").Value2 = -30309356
$ws.Range("H134").Value2 = 6287.923
$ws.Range("I134").Value2 = 2089.8667
$ws.Range("K134").Value2 = 6269.6001
$ws.Range("M134").Value2 = -3734.6001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 8779597
$ws.Range("I31").Value2 = 2798.889
$ws.Range("J31").Value2 = 16678715
$ws.Range("K31").Value2 = 2798.889
$ws.Range("L31").Value2 = 16678715
$ws.Range("M31").Value2 = -2503.889
$ws.Range("N31").Value2 = -16679305
$ws.Range("H34").Value2 = 8779597
$ws.Range("I34").Value2 = 2798.889
$ws.Range("J34").Value2 = 16678715
$ws.Range("K34").Value2 = 2798.889
$ws.Range("L34").Value2 = 16678715
$ws.Range("M34").Value2 = -2596.889
$ws.Range("N34").Value2 = -16679119
$ws.Range("H41").Value2 = 10704.667
$ws.Range("I41").Value2 = 10704.667
$ws.Range("K41").Value2 = 10704.667
$ws.Range("M41").Value2 = -10276.667
$ws.Range("H58").Value2 = 14712699
$ws.Range("I58").Value2 = 33335056
$ws.Range("J58").Value2 = 10837.368
$ws.Range("K58").Value2 = 33335056
$ws.Range("L58").Value2 = 10837.368
$ws.Range("M58").Value2 = -33334853
$ws.Range("N58").Value2 = -11243.368
$ws.Range("H62").Value2 = 6217.25
$ws.Range("I62").Value2 = 7571.25
$ws.Range("J62").Value2 = 4863.25
$ws.Range("K62").Value2 = 7571.25
$ws.Range("L62").Value2 = 4863.25
$ws.Range("M62").Value2 = -6947.25
$ws.Range("N62").Value2 = -6111.25
$ws.Range("H65").Value2 = 6217.25
$ws.Range("I65").Value2 = 7571.25
$ws.Range("J65").Value2 = 4863.25
$ws.Range("K65").Value2 = 37856.25
$ws.Range("L65").Value2 = 24316.25
$ws.Range("M65").Value2 = -34736.25
$ws.Range("N65").Value2 = -30556.25
$ws.Range("H86").Value2 = 7005779.5
$ws.Range("I86").Value2 = 15631502
$ws.Range("K86").Value2 = 15631502
$ws.Range("M86").Value2 = -15630379
$ws.Range("H89").Value2 = 7005779.5
$ws.Range("I89").Value2 = 15631502
$ws.Range("K89").Value2 = 78157510
$ws.Range("M89").Value2 = -78151894
$ws.Range("H105").Value2 = 11911917
$ws.Range("I105").Value2 = 35714904
$ws.Range("J105").Value2 = 10423.5
$ws.Range("K105").Value2 = 35714904
$ws.Range("L105").Value2 = 10423.5
$ws.Range("M105").Value2 = -35713157
$ws.Range("N105").Value2 = -13917.5
$ws.Range("H134").Value2 = 10373.723
$ws.Range("I134").Value2 = 2150
$ws.Range("K134").Value2 = 6450
$ws.Range("M134").Value2 = -3915
$ws.Range("H136").Value2 = 14712699
$ws.Range("I136").Value2 = 33335056
$ws.Range("J136").Value2 = 10837.368
$ws.Range("K136").Value2 = 100005168
$ws.Range("L136").Value2 = 32512.104
$ws.Range("M136").Value2 = -100002618
$ws.Range("N136").Value2 = -37612.104

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value2 = 69.77778000000001
$ws.Range("I40").Value2 = 58
$ws.Range("J40").Value2 = 84.5
$ws.Range("K40").Value2 = 232
$ws.Range("L40").Value2 = 338
$ws.Range("M40").Value2 = -163
$ws.Range("N40").Value2 = -476
$ws.Range("H68").Value2 = 25003946
$ws.Range("J68").Value2 = 66674470
$ws.Range("L68").Value2 = 200023410
$ws.Range("N68").Value2 = -200025032
$ws.Range("H71").Value2 = 25003946
$ws.Range("J71").Value2 = 66674470
$ws.Range("L71").Value2 = 600070230
$ws.Range("N71").Value2 = -600078342
$ws.Range("H112").Value2 = 1812.2858
$ws.Range("I112").Value2 = 1812.2858
$ws.Range("J112").Value2 = 0
$ws.Range("K112").Value2 = 5436.857400000001
$ws.Range("L112").Value2 = 0
$ws.Range("M112").Value2 = -4328.857400000001
$ws.Range("N112").ClearContents()
$ws.Range("H121").Value2 = 1524
$ws.Range("I121").Value2 = 100
$ws.Range("J121").Value2 = 1998.6666
$ws.Range("K121").Value2 = 300
$ws.Range("L121").Value2 = 5995.9998
$ws.Range("M121").Value2 = 1010
$ws.Range("N121").Value2 = -8615.9998
$ws.Range("H132").Value2 = 9547.852000000001
$ws.Range("J132").Value2 = 14170.857
$ws.Range("L132").Value2 = 127537.713
$ws.Range("N132").Value2 = -132597.713

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value2 = 3558.111
$ws.Range("I80").Value2 = 2433
$ws.Range("J80").Value2 = 4120.6665
$ws.Range("K80").Value2 = 2433
$ws.Range("L80").Value2 = 4120.6665
$ws.Range("M80").Value2 = -1435
$ws.Range("N80").Value2 = -6116.6665
$ws.Range("H83").Value2 = 3558.111
$ws.Range("I83").Value2 = 2433
$ws.Range("J83").Value2 = 4120.6665
$ws.Range("K83").Value2 = 12165
$ws.Range("L83").Value2 = 20603.3325
$ws.Range("M83").Value2 = -7173
$ws.Range("N83").Value2 = -30587.3325
$ws.Range("H97").Value2 = 907.1081
$ws.Range("I97").Value2 = 847.65625
$ws.Range("K97").Value2 = 847.65625
$ws.Range("M97").Value2 = -351.65625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 5639.115
$ws.Range("I7").Value2 = 4271.5884
$ws.Range("K7").Value2 = 4271.5884
$ws.Range("M7").Value2 = -4159.5884
$ws.Range("H55").Value2 = 43478676
$ws.Range("I55").Value2 = 100000100
$ws.Range("J55").Value2 = 661.1539
$ws.Range("K55").Value2 = 100000100
$ws.Range("L55").Value2 = 661.1539
$ws.Range("M55").Value2 = -99999927
$ws.Range("N55").Value2 = -1007.1539
$ws.Range("H82").Value2 = 3359.75
$ws.Range("I82").Value2 = 1000
$ws.Range("J82").Value2 = 4146.3335
$ws.Range("K82").Value2 = 1000
$ws.Range("L82").Value2 = 4146.3335
$ws.Range("M82").Value2 = -639
$ws.Range("N82").Value2 = -4868.3335
$ws.Range("H85").Value2 = 3359.75
$ws.Range("I85").Value2 = 1000
$ws.Range("J85").Value2 = 4146.3335
$ws.Range("K85").Value2 = 1000
$ws.Range("L85").Value2 = 4146.3335
$ws.Range("M85").Value2 = 248
$ws.Range("N85").Value2 = -6642.3335
$ws.Range("H100").Value2 = 4675.0835
$ws.Range("I100").Value2 = 3516.1667
$ws.Range("K100").Value2 = 3516.1667
$ws.Range("M100").Value2 = -2975.1667
$ws.Range("H126").Value2 = 5639.115
$ws.Range("I126").Value2 = 4271.5884
$ws.Range("K126").Value2 = 12814.7652
$ws.Range("M126").Value2 = -10344.7652
$ws.Range("H132").Value2 = 9096792
$ws.Range("I132").Value2 = 16131567
$ws.Range("K132").Value2 = 48394701
$ws.Range("M132").Value2 = -48392171
$ws.Range("H136").Value2 = 10678.192
$ws.Range("J136").Value2 = 13473.685
$ws.Range("L136").Value2 = 40421.055
$ws.Range("N136").Value2 = -45521.055

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value2 = 852.46155
$ws.Range("I107").Value2 = 538.1
$ws.Range("J107").Value2 = 1900.3334
$ws.Range("K107").Value2 = 1614.3
$ws.Range("L107").Value2 = 5701.0002
$ws.Range("M107").Value2 = 305.6999999999998
$ws.Range("N107").Value2 = -9541.0002
$ws.Range("H126").Value2 = 3209
$ws.Range("I126").Value2 = 951.125
$ws.Range("J126").Value2 = 7724.75
$ws.Range("K126").Value2 = 2853.375
$ws.Range("L126").Value2 = 23174.25
$ws.Range("M126").Value2 = -383.375
$ws.Range("N126").Value2 = -28114.25

"Updated 249 cells, cleared 2 cells across 8 sheets."
